# Data Formatting and Exploration Finished
# Rebuilds the single "Sheet1" crime-index sheet into three sheets:
#   2017-2013 (existing data, renamed, re-ordered to front)
#   2012-2011 (new sheet, subset/variant of the index list)
#   2020-2018 (copy of the 2017-2013 layout + an extra "Unknown" row)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename the existing sheet to "2017-2013" - content/styling unchanged.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "2017-2013"

# ---------------------------------------------------------------------
# 2) Build "2012-2011" by duplicating "2017-2013" (so it inherits the
#    already best-fit column) right after it, then stripping the cell
#    styling that "2017-2013" had (the authored sheet uses the default
#    Normal style throughout, with no per-cell "s" attribute).
# ---------------------------------------------------------------------
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "2012-2011"
$ws2.Cells.ClearFormats()

# ---------------------------------------------------------------------
# 3) Duplicate "2017-2013" again (keeps its cell styling) to build
#    "2020-2018", placed after "2012-2011" so the final left-to-right
#    order is 2017-2013, 2012-2011, 2020-2018 and sheetId continues
#    1,2,3.
# ---------------------------------------------------------------------
$ws1.Copy([System.Reflection.Missing]::Value, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "2020-2018"
$ws3.Range("A27").Value = "Unknown"
$ws3.Range("B42").Select()

# ---------------------------------------------------------------------
# 4) Populate "2012-2011" with its own ordered list of crime categories
#    (two brand new labels: "Opposition to cheque laws" and
#    "Crimes related with firearms and knifes"), then drop the four
#    trailing rows it doesn't need (27 -> 23 rows).
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "Total"
$ws2.Range("A2").Value = "Homicide"
$ws2.Range("A3").Value = "Assault"
$ws2.Range("A4").Value = "Sexual crimes"
$ws2.Range("A5").Value = "Kidnapping"
$ws2.Range("A6").Value = "Defamation"
$ws2.Range("A7").Value = "Theft"
$ws2.Range("A8").Value = "Robbery"
$ws2.Range("A9").Value = "Opposition to cheque laws"
$ws2.Range("A10").Value = "Swindling"
$ws2.Range("A11").Value = "Production and commerce of drugs"
$ws2.Range("A12").Value = "Use and purchase of drugs"
$ws2.Range("A13").Value = "Forgery"
$ws2.Range("A14").Value = "Bad treatment"
$ws2.Range("A15").Value = "Embezzlement"
$ws2.Range("A16").Value = "Bribery"
$ws2.Range("A17").Value = "Smuggling"
$ws2.Range("A18").Value = "Traffic crimes"
$ws2.Range("A19").Value = "Forestry crimes"
$ws2.Range("A20").Value = "Crimes related with firearms and knifes"
$ws2.Range("A21").Value = "Opposition  to the  Bankruptcy  and Enforcement Law"
$ws2.Range("A22").Value = "Opposition to the Military Criminal Law"
$ws2.Range("A23").Value = "Other crimes"
$ws2.Range("A24:A27").EntireRow.Delete()

# ---------------------------------------------------------------------
# 5) Column widths: let Excel compute "best fit" on the sheets whose
#    content actually changed (stamps the bestFit flag like the
#    authored file); "2012-2011" keeps the best-fit column it inherited
#    from the copy since its longest label is unchanged.
# ---------------------------------------------------------------------
$ws1.UsedRange.Columns.AutoFit()
$ws3.UsedRange.Columns.AutoFit()

# ---------------------------------------------------------------------
# 6) View state: "2012-2011" ends up the active/selected tab, with A7
#    selected; "2017-2013" gets the whole of column A selected.
# ---------------------------------------------------------------------
$ws1.Range("A1:A1048576").Select()
$ws2.Activate()
$ws2.Range("A7").Select()
